# chore: adapt column header formatting to respective input file names
#
# Rename header row cells from the generic "<Name>_old" / "<Name>_new"
# suffixes to the concrete format-version suffixes "<Name>_FV2210" and
# "<Name>_FV2304" respectively, then turn the data range into a real Excel
# Table (so there's a header-row autofilter) and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$lastCol = $usedRange.Columns.Count

for ($col = 1; $col -le $lastCol; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $header = $cell.Value2
    if ($header -like "*_old") {
        $cell.Value = ($header -replace '_old$', '_FV2210')
    } elseif ($header -like "*_new") {
        $cell.Value = ($header -replace '_new$', '_FV2304')
    }
}

# Build the A1:U61 range reference for the header+data block.
$lastColLetter = $ws.Cells.Item(1, $lastCol).Address($false, $false) -replace '\d+$', ''
$tableRange = $ws.Range("A1:" + $lastColLetter + $lastRow)

# Convert the range into a native Excel table ("Table1") with a header row
# / autofilter, matching the column names we just renamed above.
$tbl = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$tbl.Name = "Table1"

# Freeze the header row: select the first cell below the header and freeze.
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
